$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 85) with the July 2025 recurrence metrics.
$row = 85
$ws.Cells.Item($row, 1).Value = "2025-07"
$ws.Cells.Item($row, 2).Value = 1
$ws.Cells.Item($row, 3).Value = 207
$ws.Cells.Item($row, 4).Value = 0.4830917874396135
